$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row
$ws.Range("B1").Value = "Densidade"

# Update data rows for remaining columns
$ws.Range("B2").Value = 1000
$ws.Range("B3").Value = 2000

# Clear columns C and D entirely (remove old P0X/C helper columns data)
$ws.Range("C1:D9").Clear()

# Remove now-unused rows 4 through 9
$ws.Range("A4:D9").Clear()
